$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data block (dates serial, year, y0_forecast, year+1, y1_forecast)
$data = @(
    @(39400, 2007, 2.070003986395053, 2008, -0.3549868696899106),
    @(39765, 2008, 0.517569958955022, 2009, -5.168396053267498),
    @(40130, 2009, -3.956152295564896, 2010, -1.314964327391877),
    @(40494, 2010, 1.234995474941392, 2011, 1.001424185348321),
    @(40862, 2011, 0.899360810820804, 2012, 0.475544341751033),
    @(41228, 2012, 0.9010266119894084, 2013, 1.506358095275817),
    @(41592, 2013, 0.02019328874804938, 2014, -1.194807813319188),
    @(41957, 2014, 0.1729981757035093, 2015, 0.6265079396372775),
    @(42321, 2015, 0.09752710595589686, 2016, -1.022506370243093),
    @(42689, 2016, -0.5280591151586633, 2017, -0.7240982069264934),
    @(43053, 2017, 0.07201851318385799, 2018, 1.255028673974046),
    @(43418, 2018, 0.3727661260635617, 2019, -3.305525567352929),
    @(43783, 2019, -0.801759526476209, 2020, 1.431264289671219),
    @(44159, 2020, -1.103489789942047, 2021, 2.294626310579817),
    @(44525, 2021, 0.9704846793491928, 2022, -0.902682013141165),
    @(44890, 2022, -0.7009264669202708, 2023, 1.220523709718857),
    @(45254, 2023, 0.3928252664241905, 2024, 0.4517021897791018),
    @(45618, 2024, 0.3224026462283813, 2025, -0.7185940249203049)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $row++
}

# Ensure the newly added row (19) uses the same formatting as the other date cells in column A
$ws.Range("A18").Copy()
$ws.Range("A19").PasteSpecial(-4122)
